$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 838 (2026/12/29 block). All rows
# from 838 downward shift down by one; Excel auto-extends the used range.
$ws.Rows.Item(838).Insert()

# New row 838: 2026/02/23 (Mon), time 5, rank 21.
# Force column A to be stored as literal text (not auto-parsed into a date
# serial number) the same way the rest of the "日付" column is stored.
$ws.Range("A838").NumberFormat = "@"
$ws.Range("A838").Value = "2026/02/23"
$ws.Range("A838").Style = "Normal"

$ws.Range("B838").Value = "月"
$ws.Range("C838").Value = 5
$ws.Range("D838").Value = 21
